$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value = 756
$wsExhibition.Range("F12").Value = 220
$wsExhibition.Range("F15").Value = 121
$wsExhibition.Range("F17").Value = 512
$wsExhibition.Range("F18").Value = 8030
$wsExhibition.Range("F19").Value = 618

# Sheet "本地生活" (Local Life) update
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 5569

# Sheet "全部类型" (All Types) updates
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F3").Value = 5569
$wsAllTypes.Range("F13").Value = 756
$wsAllTypes.Range("F20").Value = 220
$wsAllTypes.Range("F25").Value = 121
$wsAllTypes.Range("F29").Value = 512
$wsAllTypes.Range("F30").Value = 8030
$wsAllTypes.Range("F33").Value = 618
